# Örnek 10 - Metin İşlemleri.xlsx
# Fill in the "metin işlemleri" (text operations) example row 5 with formulas
# that combine A1/B1, and compute length / upper / lower / mid / value,
# add the currency number format to the VALUE() result, fill in the
# student info block (F8:F10), and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: text-function demo (A5:C5 is merged; result cells D5..H5) ---
$ws.Range("A5").Formula = '=A1&" "&B1'
$ws.Range("D5").Formula = "=LEN(A5)"
$ws.Range("E5").Formula = "=UPPER(A5)"
$ws.Range("F5").Formula = "=LOWER(A5)"
$ws.Range("G5").Formula = "=MID(A5,5,3)"
$ws.Range("H5").Formula = "=VALUE(D5)"

# Give E5 a wrap-text style, then propagate that same visual style down to
# the (still empty) E6 cell below it.
$ws.Range("E5").WrapText = $true
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# H5 holds a numeric result (the character count, as a number) -- format it
# as Turkish Lira currency, same as the "lira" column header describes.
$ws.Range("H5").NumberFormat = '_-[$₺-41F]* #,##0.00_-;\-[$₺-41F]* #,##0.00_-;_-[$₺-41F]* "-"??_-;_-@_-'

# Row 6 grew a bit taller in the edited workbook.
$ws.Rows("6:6").RowHeight = 20.25

# --- Student info block ---
$ws.Range("F8").Value = 20215070019
$ws.Range("F9").Value = "KÜBRA ÇABUK"
$ws.Range("F10").Value = "YBS"

# Minor column width tweaks.
$ws.Columns("D:D").ColumnWidth = 28.43
$ws.Columns("E:E").ColumnWidth = 23.29

# Move the active selection to G5, matching where the editor left off.
[void]$ws.Range("G5").Select()
